$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.815.34'
$ws.Range("E2").Value = '  +4.01%  '

$ws.Range("D3").Value = '3.251.47'
$ws.Range("E3").Value = '  +2.12%  '

$ws.Range("E5").Value = '  -1.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.92'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("E7").Value = '  +5.32%  '

$ws.Range("D8").Value = '3.248.00'
$ws.Range("E8").Value = '  +2.10%  '

$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("E10").Value = '  +1.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '39.30'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0970'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +9.50%  '

$ws.Range("E13").Value = '  +2.41%  '

$ws.Range("D14").Value = '3.765.23'
$ws.Range("E14").Value = '  +2.46%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.33'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +3.61%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.13'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +0.53%  '

$ws.Range("D17").Value = '3.253.82'
$ws.Range("E17").Value = '  +2.18%  '

$ws.Range("E18").Value = '  -3.44%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.70'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +1.73%  '

$ws.Range("D20").Value = '56.749.91'
$ws.Range("E20").Value = '  +4.13%  '

$ws.Range("E21").Value = '  +1.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000109'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +8.91%  '

$ws.Range("E23").Value = '  +0.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '295.39'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +7.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.19'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +2.50%  '

$ws.Range("E26").Value = '  -3.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.18'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +1.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.65'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -4.98%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.27'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -2.48%  '

$ws.Range("E31").Value = '  -1.04%  '

$ws.Range("E32").Value = '  +0.04%  '

$ws.Range("E33").Value = '  +2.21%  '

$ws.Range("E34").Value = '  -4.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '39.80'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +7.29%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0485'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -5.01%  '

$ws.Range("E37").Value = '  +2.25%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '51.36'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +0.94%  '

$ws.Range("E39").Value = '  -0.01%  '

$ws.Range("E40").Value = '  -4.97%  '

$ws.Range("E41").Value = '  +0.99%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '136.74'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +4.70%  '

$ws.Range("E43").Value = '  +3.48%  '

$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.90'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -2.49%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.97'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -4.77%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.02'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -1.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.281'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -3.47%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.30'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -0.09%  '

$ws.Range("E49").Value = '  +3.10%  '

$ws.Range("D50").Value = '2.156.59'
$ws.Range("E50").Value = '  +3.10%  '

$ws.Range("E51").Value = '  -5.57%  '
